$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column values are written as plain text (matches original inlineStr cells),
# since several new values (e.g. "211.06") would otherwise be auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.318.65"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.566.51"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D5").Value = "211.06"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "0.491"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "44.34"
$ws.Range("E8").Value = "  -4.18%  "
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "1.790.51"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "1.573.38"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "3.67"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "28.309.66"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "0.513"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "61.03"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "227.49"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "0.0₃0679"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "3.94"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "8.94"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "0.104"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "6.33"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("D35").Value = "1.381.18"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("D40").Value = "0.0163"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "0.521"
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("D42").Value = "1.93"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "0.782"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "5.33"
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("D47").Value = "62.23"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "0.915"
$ws.Range("E48").Value = "  -6.42%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.702.97"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "85.62"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.89%  "

# Restore the default (unstyled) look for column D so no stray number-format style sticks around.
$ws.Range("D2:D51").Style = "Normal"

